$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a "last changed" date serial for every data
# row (2-480). This automated update bumps that date by one day
# (45179 -> 45180, i.e. 2023-09-10 -> 2023-09-11) for every row.
$ws.Range("C2:C480").Value = 45180
